# Fixing typo in starter letter:
# "...feedback on IGB functionality or suggesting ways we can improve our work."
# becomes
# "...feedback on IGB functionality."

$d = $word.ActiveDocument

# Locate the exact span of text that needs to change, scoped to a Range so we
# don't disturb the rest of the paragraph (e.g. the preceding "feedback" run).
$oldText = " on IGB functionality or suggesting ways we can improve our work. "
$finder = $d.Content
$found = $finder.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text to fix typo in starter letter."
}
$targetRange = $d.Range($finder.Start, $finder.End)

# Replace the run's contents with equivalent WordprocessingML that keeps the
# original run (now reading " on IGB ") and appends two new same-formatted
# runs for "functionality" and ". " -- exactly how Word splits a run when
# the trailing part of its text is edited out.
$replacementXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r w:rsidR="00620B28">
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:sz w:val="22"/>
                <w:highlight w:val="cyan"/>
              </w:rPr>
              <w:t xml:space="preserve"> on IGB </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:sz w:val="22"/>
                <w:highlight w:val="cyan"/>
              </w:rPr>
              <w:t>functionality</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:sz w:val="22"/>
                <w:highlight w:val="cyan"/>
              </w:rPr>
              <w:t xml:space="preserve">. </w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$targetRange.InsertXML($replacementXml)
